$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H: new "学费" (tuition fee) column -----------------------------

# Header text
$ws.Range("H1").Value = "学费"

# Fee values per contiguous block of rows (matches the major / study-mode
# groupings already present in columns A-D)
$ranges = @(
    @{ First = 2;   Last = 23;  Fee = 10000 },
    @{ First = 24;  Last = 56;  Fee = 26000 },
    @{ First = 57;  Last = 57;  Fee = 30000 },
    @{ First = 58;  Last = 75;  Fee = 29000 },
    @{ First = 76;  Last = 105; Fee = 35000 },
    @{ First = 106; Last = 112; Fee = 10000 },
    @{ First = 113; Last = 142; Fee = 15000 }
)

foreach ($r in $ranges) {
    $rng = $ws.Range("H" + $r.First + ":H" + $r.Last)
    $rng.Value = $r.Fee
}

# Copy formatting (font/border/alignment) from existing cells instead of
# rebuilding it property-by-property, so the new column reuses the exact
# same font/border objects as the rest of the sheet.

# Body cells -> match the other data columns (e.g. G2's style)
[void]$ws.Range("G2").Copy()
[void]$ws.Range("H2:H142").PasteSpecial(-4122)

# Header cell -> match the other header cells (e.g. A1's style)
[void]$ws.Range("A1").Copy()
[void]$ws.Range("H1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Sheet view state --------------------------------------------------
$excel.ActiveWindow.ScrollRow = 95
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E96").Select()
